$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.09934133333334
$ws.Range("H2").Value = 63.29802400000001
$ws.Range("I2").Value = 0.2917236204149438
$ws.Range("J2").Value = 0.2917236204149438
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.842817666666666
$ws.Range("N2").Value = 14.528453
$ws.Range("O2").Value = 0.08128949930032948
$ws.Range("P2").Value = 0.0812894993003295
$ws.Range("Q2").Value = 102.1802629640969
$ws.Range("R2").Value = 919.622366676872
$ws.Range("S2").Value = 0.02371406703761016
$ws.Range("T2").Value = 0.02371406703761016

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.09934133333334
$ws.Range("H3").Value = 63.29802400000001
$ws.Range("I3").Value = 0.2917236204149438
$ws.Range("J3").Value = 0.2917236204149438
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 38.12230933333333
$ws.Range("N3").Value = 114.366928
$ws.Range("O3").Value = 0.6399050410691925
$ws.Range("P3").Value = 0.6399050410691927
$ws.Range("Q3").Value = 804.3556170389191
$ws.Range("R3").Value = 7239.200553350272
$ws.Range("S3").Value = 0.1866754153024781
$ws.Range("T3").Value = 0.1866754153024782

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.09934133333334
$ws.Range("H4").Value = 63.29802400000001
$ws.Range("I4").Value = 0.2917236204149438
$ws.Range("J4").Value = 0.2917236204149438
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.60982066666667
$ws.Range("N4").Value = 49.82946200000001
$ws.Range("O4").Value = 0.2788054596304779
$ws.Range("P4").Value = 0.2788054596304779
$ws.Range("Q4").Value = 350.4562757314543
$ws.Range("R4").Value = 3154.106481583089
$ws.Range("S4").Value = 0.08133413807485546
$ws.Range("T4").Value = 0.08133413807485548

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 35.81943766666667
$ws.Range("H5").Value = 107.458313
$ws.Range("I5").Value = 0.4952465516465762
$ws.Range("J5").Value = 0.4952465516465762
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.842817666666666
$ws.Range("N5").Value = 14.528453
$ws.Range("O5").Value = 0.08128949930032948
$ws.Range("P5").Value = 0.0812894993003295
$ws.Range("Q5").Value = 173.4670055421988
$ws.Range("R5").Value = 1561.203049879789
$ws.Range("S5").Value = 0.04025834421356494
$ws.Range("T5").Value = 0.04025834421356495

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 35.81943766666667
$ws.Range("H6").Value = 107.458313
$ws.Range("I6").Value = 0.4952465516465762
$ws.Range("J6").Value = 0.4952465516465762
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 38.12230933333333
$ws.Range("N6").Value = 114.366928
$ws.Range("O6").Value = 0.6399050410691925
$ws.Range("P6").Value = 0.6399050410691927
$ws.Range("Q6").Value = 1365.519682874718
$ws.Range("R6").Value = 12289.67714587247
$ws.Range("S6").Value = 0.3169107649707783
$ws.Range("T6").Value = 0.3169107649707784

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 35.81943766666667
$ws.Range("H7").Value = 107.458313
$ws.Range("I7").Value = 0.4952465516465762
$ws.Range("J7").Value = 0.4952465516465762
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.60982066666667
$ws.Range("N7").Value = 49.82946200000001
$ws.Range("O7").Value = 0.2788054596304779
$ws.Range("P7").Value = 0.2788054596304779
$ws.Range("Q7").Value = 594.9544360241784
$ws.Range("R7").Value = 5354.589924217606
$ws.Range("S7").Value = 0.1380774424622329
$ws.Range("T7").Value = 0.1380774424622329

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 15.40769666666667
$ws.Range("H8").Value = 46.22309
$ws.Range("I8").Value = 0.2130298279384801
$ws.Range("J8").Value = 0.2130298279384801
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.842817666666666
$ws.Range("N8").Value = 14.528453
$ws.Range("O8").Value = 0.08128949930032948
$ws.Range("P8").Value = 0.0812894993003295
$ws.Range("Q8").Value = 74.61666561997443
$ws.Range("R8").Value = 671.5499905797699
$ws.Range("S8").Value = 0.01731708804915439
$ws.Range("T8").Value = 0.01731708804915439

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 15.40769666666667
$ws.Range("H9").Value = 46.22309
$ws.Range("I9").Value = 0.2130298279384801
$ws.Range("J9").Value = 0.2130298279384801
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 38.12230933333333
$ws.Range("N9").Value = 114.366928
$ws.Range("O9").Value = 0.6399050410691925
$ws.Range("P9").Value = 0.6399050410691927
$ws.Range("Q9").Value = 587.3769784408355
$ws.Range("R9").Value = 5286.39280596752
$ws.Range("S9").Value = 0.1363188607959361
$ws.Range("T9").Value = 0.1363188607959362

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 15.40769666666667
$ws.Range("H10").Value = 46.22309
$ws.Range("I10").Value = 0.2130298279384801
$ws.Range("J10").Value = 0.2130298279384801
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.60982066666667
$ws.Range("N10").Value = 49.82946200000001
$ws.Range("O10").Value = 0.2788054596304779
$ws.Range("P10").Value = 0.2788054596304779
$ws.Range("Q10").Value = 255.9190785197311
$ws.Range("R10").Value = 2303.27170667758
$ws.Range("S10").Value = 0.05939387909338956
$ws.Range("T10").Value = 0.05939387909338957
